$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.003.33'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -5.13%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.300.38'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -5.60%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '561.78'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -4.60%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '126.54'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -5.32%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.296.48'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -5.66%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.477'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.63%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.32'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -4.79%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.373'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -3.33%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.870.45'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -5.50%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.119'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.56%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.309.22'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -5.44%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000167'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -6.85%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '24.53'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.25%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '61.073.79'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -4.93%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.40'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.93%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.62'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.40%  '

$ws.Range("E21").Value = '  -10.21%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '351.25'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -9.00%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.553'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -4.36%  '

$ws.Range("E24").Value = '  +0.10%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.435.66'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -5.52%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '69.04'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -7.11%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000106'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -7.70%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.30%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.11'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.72%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.82'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -3.75%  '

$ws.Range("E31").Value = '  -6.12%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.09'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -6.99%  '

$ws.Range("E33").Value = '  -0.06%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.148'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -4.09%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.330.55'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -5.52%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '22.46'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -3.39%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.18'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -4.26%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.74'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.09%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '160.18'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -2.62%  '

$ws.Range("E40").Value = '  -4.62%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0754'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -3.87%  '

$ws.Range("E42").Value = '  +0.06%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.95'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.20%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.32'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.85%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.740'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -8.20%  '

$ws.Range("E46").Value = '  -6.37%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.55'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -5.77%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.08'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -9.58%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.66'
$ws.Range("D49").ClearFormats()

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.861'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -6.37%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '20.78'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.68%  '
